$d = $word.ActiveDocument

# --- 1. Insert a new "Acknowledgments" Heading1 paragraph (with bookmark)
#        right before the existing "Credits" paragraph. --------------------
$creditsIndex = $d.Paragraphs.Count
$creditsPara = $d.Paragraphs($creditsIndex)
$creditsPara.Range.InsertParagraphBefore()

$newPara = $d.Paragraphs($creditsIndex)
$newPara.Style = "Heading1"
$newPara.Range.Text = "Acknowledgments"

$bmStart = $newPara.Range.Start
$bmEnd = $newPara.Range.End
$d.Bookmarks.Add("acknowledgments", $d.Range($bmStart, $bmEnd))

# --- 2. Replace the Credits placeholder text. ------------------------------
$d.Content.Find.Execute("Credits placeholder.", $true, $false, $false, $false, $false, $true, 1, $false, "Some materials included in this export came from the following casebooks.", 2)

# --- 3. Give the (previously empty) final section explicit page setup. ----
$ps = $d.PageSetup
$ps.PageWidth = 612
$ps.PageHeight = 792
$ps.TopMargin = 72
$ps.BottomMargin = 72
$ps.LeftMargin = 72
$ps.RightMargin = 72
$ps.HeaderDistance = 36
$ps.FooterDistance = 36
$ps.Gutter = 0
$ps.TextColumns.Spacing = 36

Write-Output "ok"
